$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Elimna EC anteriores y se agregan nuevos" -- the previous-period account
# statements (Periodo Mora / Valor Mora rows 16-19) are replaced: the macro's
# refreshed export reverses the period ordering (1805,1806,1807,1903 ->
# 1903,1807,1806,1805) and the "Valor Mora" (column F) follows the row that
# now holds period "1903".
$ws.Range("E16").Value = "1903"
$ws.Range("F16").Value = 33125

$ws.Range("E17").Value = "1807"
$ws.Range("F17").Value = 31249

$ws.Range("E18").Value = "1806"
$ws.Range("F18").Value = 31249

$ws.Range("E19").Value = "1805"
$ws.Range("F19").Value = 31249
